# Apply updated Price (column D) values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$priceUpdates = @{
    2 = 9.470000000000001
    3 = 8.68
    4 = 13.29
    5 = 16.93
    6 = 2.99
    7 = 9.07
    8 = 7.47
    9 = 5.51
    10 = 24.51
    11 = 0.79
    12 = 1
    13 = 5.93
    14 = 1.3
    15 = 6.52
    16 = 3.71
    17 = 0.9399999999999999
    18 = 0.71
    19 = 0.93
    20 = 15.69
    21 = 3.97
    22 = 8.130000000000001
    23 = 4.56
    24 = 1.08
    25 = 6.66
    26 = 2.95
    27 = 8.359999999999999
    29 = 4.96
    30 = 0.8100000000000001
    31 = 4.21
    32 = 0.26
    33 = 1.05
    34 = 1.69
    35 = 1.53
    36 = 2.83
    37 = 9.93
    39 = 1.65
    40 = 5.07
    41 = 26.16
    42 = 0.64
    43 = 1.31
    44 = 2.03
    46 = 3.87
    47 = 1.8
    48 = 1.09
    49 = 3.41
    50 = 2.91
    51 = 6.44
    52 = 3.28
    53 = 1.86
    54 = 1.52
    55 = 6.02
    56 = 1.58
    57 = 3.28
    58 = 4.38
    59 = 1.33
    60 = 10.73
    61 = 8.16
    62 = 17.87
    63 = 1.12
    64 = 9.710000000000001
    65 = 9.07
    66 = 4.3
    67 = 2.4
    68 = 5.21
    69 = 3.98
    70 = 1.1
    71 = 0.57
    72 = 1.42
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $priceUpdates[$row]
}

